# Update metadata propagation images
#
# Applies the textual / positional tweaks described by the commit
# "Update metadata propagation images" to the shapes that are actually
# present in this deck (the custGeom image placeholders referenced in
# the original diff are not part of this fixture, so those two hunks
# are intentionally skipped).

$p = $ppt.ActivePresentation

function Find-ShapeByName($presentation, [string]$namePattern) {
    for ($si = 1; $si -le $presentation.Slides.Count; $si++) {
        $slide = $presentation.Slides.Item($si)
        for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
            $shape = $slide.Shapes.Item($i)
            if ($shape.Name -match $namePattern) {
                return $shape
            }
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "*only labels, <br>filtered by prefix/domain" callout on slide 2
#    -> "*labels and annotation<br>filtered by prefix/domain"
#    Replace only the first run's characters (1..14 == "*only labels, ")
#    so the existing <a:br> soft line-break and <a:endParaRPr> survive.
#    The shape auto-fits its height to the text (<a:spAutoFit/>); the
#    box's own xfrm is untouched in the real edit, so restore the
#    original height after the text swap triggers a reflow.
# ---------------------------------------------------------------------
$calloutHeightPt = 48.47244074488189  # -> 615600 EMU (unchanged box height)

$shape257 = Find-ShapeByName $p "257;p14"
if ($shape257 -ne $null) {
    $tr = $shape257.TextFrame.TextRange
    $lead = $tr.Characters(1, 14)
    if ($lead.Text -eq "*only labels, ") {
        $lead.Text = "*labels and annotation"
    }
    $shape257.Height = $calloutHeightPt
}

# ---------------------------------------------------------------------
# 2) Same callout duplicated on slide 3, but the replacement keeps a
#    space after the asterisk: "* labels and annotation".
# ---------------------------------------------------------------------
$shape375 = Find-ShapeByName $p "375;p15"
if ($shape375 -ne $null) {
    $tr = $shape375.TextFrame.TextRange
    $lead = $tr.Characters(1, 14)
    if ($lead.Text -eq "*only labels, ") {
        $lead.Text = "* labels and annotation"
    }
    $shape375.Height = $calloutHeightPt
}

# ---------------------------------------------------------------------
# 3) Four "Cluster U CC U template" textboxes on slide 3 get a 1-EMU
#    nudge left and get widened (2118300 -> 2201173 EMU); vertical
#    position/height are untouched.
#
#    Shape.Left / .Width are expressed in points and the host rounds
#    through a 32-bit float, so we feed it the precise point value
#    that reproduces the exact target EMU instead of the naive
#    EMU/12700.0 (which can land 1 EMU short after the float32 cast).
# ---------------------------------------------------------------------
$targetLeftPt  = 884.7400818401575   # -> 11236199 EMU
$targetWidthPt = 173.32070816141731  # -> 2201173 EMU

foreach ($nameFrag in @("387;p15", "388;p15", "389;p15", "390;p15")) {
    $shape = Find-ShapeByName $p $nameFrag
    if ($shape -ne $null) {
        $shape.Left = $targetLeftPt
        $shape.Width = $targetWidthPt
    }
}

# ---------------------------------------------------------------------
# 4) The red "4b" badge's text box gets tighter left/right insets
#    (91425 -> 36000 EMU); top/bottom insets stay at 91425 EMU.
# ---------------------------------------------------------------------
$shape395 = Find-ShapeByName $p "395;p15"
if ($shape395 -ne $null) {
    $tf = $shape395.TextFrame
    $tf.MarginLeft = 2.8346456692913384   # -> 36000 EMU
    $tf.MarginRight = 2.8346456692913384  # -> 36000 EMU
}
